# Remove the "Egg_FoundID" lookup column from the Eggs sheet now that
# getGottenEggs() resolves found-eggs directly (and also returns the next
# eggs) instead of relying on a stored column.
$wb = $excel.ActiveWorkbook

$wsUsers     = $wb.Worksheets.Item("Users")
$wsEggs      = $wb.Worksheets.Item("Eggs")
$wsEggLinks  = $wb.Worksheets.Item("EggLinks")
$wsFoundList = $wb.Worksheets.Item("FoundList")

# Eggs!B ("Egg_FoundID") is gone; EggLocation/EggRedirect/EggValue shift left.
$wsEggs.Range("B1").EntireColumn.Delete()

# Restore each sheet's view (zoom level + selection), then finish on
# FoundList so it becomes the active / selected tab, matching the saved
# workbook state.
$wsUsers.Activate()
$excel.ActiveWindow.Zoom = 160
$wsUsers.Range("E1").Select()

$wsEggs.Activate()
$excel.ActiveWindow.Zoom = 130
$wsEggs.Range("B1:B1048576").Select()

$wsEggLinks.Activate()
$excel.ActiveWindow.Zoom = 205
$wsEggLinks.Range("A5").Select()

$wsFoundList.Activate()
$excel.ActiveWindow.Zoom = 190
$wsFoundList.Range("B3").Select()
